# Apply the commit's data refresh to the "Logs" worksheet of the Jira Status Report.
# The dispatcher run produced a new batch of log rows (later timestamps and
# new ticket numbers) while the overall table layout stayed the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: TicketCreationDate timestamps (row 2 -> row 8)
$ws.Range("A2").Value = "12-Feb-2022 06:16:35"
$ws.Range("A3").Value = "12-Feb-2022 06:16:53"
$ws.Range("A4").Value = "12-Feb-2022 06:16:57"
$ws.Range("A5").Value = "12-Feb-2022 06:17:27"
$ws.Range("A6").Value = "12-Feb-2022 06:17:44"
$ws.Range("A7").Value = "12-Feb-2022 06:17:47"
$ws.Range("A8").Value = "12-Feb-2022 06:18:04"

# Column C: Ticket Number updates for the successful rows
$ws.Range("C3").Value = "GRD-2776"
$ws.Range("C6").Value = "GRD-2777"
$ws.Range("C8").Value = "GRD-2778"

# Restore the active selection to A3 (matches the saved sheet view state)
$ws.Range("A3").Select()
